$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Corrected")

# Insert a new column C for "adductName" (pushes formula/medMz-like cols right)
$ws.Columns.Item(3).Insert()

# Header cell: bold font (matches the other headers) but reset to the
# "Normal" base style first so it does NOT inherit the centered alignment
# that the row's default style (s=1) carries.
$ws.Cells.Item(1, 3).Style = "Normal"
$ws.Cells.Item(1, 3).Value = "adductName"
$ws.Cells.Item(1, 3).Font.Bold = $true

# New column's data rows are all 0
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 3).Value = 0
}

# Match the column width used for this header elsewhere in the workbook
$ws.Columns.Item(3).ColumnWidth = 10.166666666666666
